$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("CE").Insert()
"done"
